$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the environment URLs/values from "test7" to "test18"
$ws.Range("A2").Value = "https://test18.cliotest.com/backoffice/control/main"
$ws.Range("C2").Value = "https://test18.cliotest.com/cabicentral/control/main"
$ws.Range("D2").Value = "https://test18.cliotest.com/warehouse/control/main"
$ws.Range("F2").Value = "virtual_cabitest18"
$ws.Range("G2").Value = "test18"
$ws.Range("K2").Value = "test18"

# Update the active selection shown in the sheet view
$ws.Range("C12").Select()
